$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07496433333333334
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08378199999999998
$ws.Range("N2").Value = 0.251346
$ws.Range("O2").Value = 0.007571394704126512
$ws.Range("P2").Value = 0.007571394704126512
$ws.Range("Q2").Value = 0.006280661775333333
$ws.Range("R2").Value = 0.056525955978
$ws.Range("S2").Value = 0.007571394704126512
$ws.Range("T2").Value = 0.007571394704126512

# Row 3
$ws.Range("G3").Value = 0.07496433333333334
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("O3").Value = 0.006644889460697858
$ws.Range("P3").Value = 0.006644889460697857
$ws.Range("Q3").Value = 0.00551210244188889
$ws.Range("R3").Value = 0.049608921977
$ws.Range("S3").Value = 0.006644889460697858
$ws.Range("T3").Value = 0.006644889460697857

# Row 4
$ws.Range("G4").Value = 0.07496433333333334
$ws.Range("O4").Value = 0.9857837158351757
$ws.Range("P4").Value = 0.9857837158351755
$ws.Range("Q4").Value = 0.817732312835
$ws.Range("R4").Value = 7.359590815515
$ws.Range("S4").Value = 0.9857837158351757
$ws.Range("T4").Value = 0.9857837158351755
